$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old rows 37-52 (old content: GenFig1.R..BinomialDownsampling.R block + stray notes rows 51-52)
$ws.Range("A37:B52").EntireRow.Delete()

# Rebuild rows 37-54 with the new File/Function verification matrix entries
$ws.Range("A37").Value = "GenFig1_3.R"
$ws.Range("B37").Value = "Difficult to test. We see some things that convince us that it works: The .25 points ends up as expected, although calculated separately. The prediction works as expected, producing reasonably constant values for the predictions. The rest of the code is rather trivial."
$ws.Range("A37").Font.Bold = $true
$ws.Range("A37").HorizontalAlignment = -4131
$ws.Rows.Item(37).RowHeight = 45

$ws.Range("A38").Value = "GenFig1_3Data.R"
$ws.Range("B38").Value = "See GenFig1.R"
$ws.Range("A38").Font.Bold = $true
$ws.Range("A38").HorizontalAlignment = -4131

$ws.Range("A39").Value = "GenFig2_S4_S5.R"
$ws.Range("B39").Value = "These are fairly simple plots with results that look as expected. No further verification is deemed necessary."
$ws.Range("A39").Font.Bold = $true
$ws.Range("A39").HorizontalAlignment = 1

$ws.Range("A40").Value = "GenFig4A-C_S23.R"
$ws.Range("B40").Value = "These plots show complicated things, but all that is calculated in underlying functions with a separate test entry in this verification matrix. The figure code in itself is fairly straight-forward. No explicit testing is deemed necessary."
$ws.Range("A40").Font.Bold = $true
$ws.Range("A40").HorizontalAlignment = 1
$ws.Rows.Item(40).RowHeight = 30

$ws.Range("A41").Value = "GenFig4A-C_S23Data.R"
$ws.Range("B41").Value = "See GenFig4A-C_S23.R"
$ws.Range("A41").Font.Bold = $true
$ws.Range("A41").HorizontalAlignment = 1

$ws.Range("A42").Value = "GenFig4DE.R"
$ws.Range("B42").Value = "The code has been reviewed, and the plot shows the expected outcome. No further testing has been done."
$ws.Range("A42").Font.Bold = $true
$ws.Range("A42").HorizontalAlignment = 1

$ws.Range("A43").Value = "GenFig4DEData.R"
$ws.Range("B43").Value = "See GenFig4.R"
$ws.Range("A43").Font.Bold = $true
$ws.Range("A43").HorizontalAlignment = 1

$ws.Range("A44").Value = "GenFig5AB.R"
$ws.Range("B44").Value = "The batch correction in itself can easily be verified graphically in the figure, it is unlikely that the improvements would come if the batch correction didn't work. However, the nearest neighbor calculations should be tested. We check a few things: 1. That the dataset source and coordinates match, i.e. have the same cell ids at the same indices (commented as test 1 in the code, only tested for uncorrected, the code for corrected is identical). The knn calculations are tested with a test case in the code (Test 2)."
$ws.Range("A44").Font.Bold = $true
$ws.Range("A44").HorizontalAlignment = 1
$ws.Rows.Item(44).RowHeight = 60

$ws.Range("A45").Value = "GenFig5C-H_S24-S25.R"
$ws.Range("B45").Value = "FIg S24 shows that the Seurat processing has worked somewhat at least. In addition, we check that we identified the clusters correctly, see Test 1 and Test 2 in the code. We test that the extracted CU per cluster for cluster 0 matches that explicitly calculated in a different way (Test 3). The rest of the code is difficult to test, but the results look as expected, no large surprises."
$ws.Range("A45").Font.Bold = $true
$ws.Range("A45").HorizontalAlignment = 1
$ws.Rows.Item(45).RowHeight = 45

$ws.Range("A46").Value = "GenFigData.R"
$ws.Range("B46").Value = "This code just calls data generation functions for the datasets, there is no need to test the code in this file."
$ws.Range("A46").Font.Bold = $true
$ws.Range("A46").HorizontalAlignment = 1

$ws.Range("A47").Value = "GenFigS1-S3.R"
$ws.Range("B47").Value = "These are fairly simple plots with results that look as expected. No further verification is deemed necessary."
$ws.Range("A47").Font.Bold = $true
$ws.Range("A47").HorizontalAlignment = 1

$ws.Range("A48").Value = "GetFigS6.R"
$ws.Range("B48").Value = "This code is trivial, it just plots the data generated by GenFigS6Data - no tests were deemed needed."
$ws.Range("A48").Font.Bold = $true
$ws.Range("A48").HorizontalAlignment = 1

$ws.Range("A49").Value = "GenFigS6Data.R"
$ws.Range("B49").Value = "The code is generally difficult to test, we mostly rely on external R packages. We do test the GC calculation function though (Test 1) and that the mean calculation across groups work for the variables (Test 2)"
$ws.Range("A49").Font.Bold = $true
$ws.Range("A49").HorizontalAlignment = 1
$ws.Rows.Item(49).RowHeight = 30

$ws.Range("A50").Value = "GenFigS7-S21.R"
$ws.Range("B50").Value = "These plots are fairly complicated, but are difficult to test. The outcome of the joint plot and the individual ds plots looks similar, which is a good sign. The scatter plots of ZTNB show less error in both scatter plots and in the loess, which looks reassuring. The code has been reviewed. No further testing is done."
$ws.Range("A50").Font.Bold = $true
$ws.Range("A50").HorizontalAlignment = 1
$ws.Rows.Item(50).RowHeight = 45

$ws.Range("A51").Value = "GenFigS7-S21Data.R"
$ws.Range("B51").Value = "This code is difficult to test and is fairly straight-forward, it calls a lot of prediction methods and saves the data. The produced end results look reasonable. The code has been reviewed. No further testing is deemed necessary."
$ws.Range("A51").Font.Bold = $true
$ws.Range("A51").HorizontalAlignment = 1
$ws.Rows.Item(51).RowHeight = 30

$ws.Range("A52").Value = "GenFigS22.R"
$ws.Range("B52").Value = "No explicit tests have been done for this code. It is mainly graphical code, although there are some lines about AUC calculations. Those follow the recommended way of using the package, and yield the expected result, so no more tests were deemed necessary."
$ws.Range("A52").Font.Bold = $true
$ws.Range("A52").HorizontalAlignment = 1
$ws.Rows.Item(52).RowHeight = 45

$ws.Range("A53").Value = "GenFigS22Data.R"
$ws.Range("B53").Value = "We test the function calcFSCM with Test 1 in the code, the function trimZeros with Test 2 in the code, and the large function evaluateCondition with Test 3 in the code. Test 3 also implicitly tests genGeneData. The code for linear interpolation was not formally tested, although it produces reasonable results. The rest of the code is not formally tested, it mainly uses the other functions, the code is not that complicated. The prediction and binomial downsampling functions were not tested, but rather validated as they do improve the classifications in the figures. Estimation of counts per cell was not explicitly tested, although it produced reasonable results. That code is only a few lines."
$ws.Range("A53").Font.Bold = $true
$ws.Range("A53").HorizontalAlignment = 1
$ws.Rows.Item(53).RowHeight = 90

$ws.Range("A54").Value = "modZTNB.R"
$ws.Range("B54").Value = "This file is copied from PreseqR and slightly modified. The performance gain vs accuracy loss is tested with TCR0004."
$ws.Range("A54").Font.Bold = $true
$ws.Range("A54").HorizontalAlignment = 1

# Update selection to match the saved cursor position (B17), and ensure no frozen/scrolled topLeftCell override
$ws.Range("B17").Select()

# Best-effort: try to restore the last-saved window geometry (may be a no-op in this host)
$wb.Windows.Item(1).Left = 1560
$wb.Windows.Item(1).Top = 1485
$wb.Windows.Item(1).Width = 25170
$wb.Windows.Item(1).Height = 14715
